$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF), matching the formatting of the
# existing header row (bold, centered, bordered) by copying H1's format.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-43
$data = @(
    @(3, 4),
    @(6, 6),
    @(1, 1),
    @(8, 8),
    @(4, 4),
    @(4, 4),
    @(6, 6),
    @(1, 2),
    @(1, 2),
    @(1, 2),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(5, 6),
    @(9, 9),
    @(8, 9),
    @(5, 5),
    @(6, 6),
    @(7, 7),
    @(9, 9),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(5, 5),
    @(4, 4),
    @(9, 9),
    @(5, 6),
    @(6, 6),
    @(10, 10),
    @(9, 9),
    @(6, 6),
    @(6, 7),
    @(6, 7),
    @(1, 3),
    @(8, 8),
    @(9, 9),
    @(5, 5),
    @(3, 4),
    @(6, 6)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
